# Regenerate the "K" column (G) values for each game row on Sheet1.
# These values are the (re-)computed strikeout counts ("K") that replace
# the previous "Strike#" derived figures, as described in the commit
# message: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$kValues = @{
    2  = 2
    3  = 3
    4  = 1
    5  = 0
    6  = 0
    7  = 2
    8  = 2
    9  = 3
    10 = 2
    11 = 2
    12 = 1
    13 = 3
    14 = 0
    15 = 2
    16 = 1
    17 = 1
    18 = 3
    19 = 0
    20 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
